$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$imgUrl = "https://sp.yimg.com/ib/th?id=OPA.yaqWrHAvuTncGA474C474&o=5&pid=21.1&w=96&h=96"

# Header row (columns A-D, then F-H for imageUrl1-3; E/originalUrl added last)
$ws.Range("A1").Value = "productName"
$ws.Range("B1").Value = "brandName"
$ws.Range("C1").Value = "ourPrice"
$ws.Range("D1").Value = "originalPrice"
$ws.Range("F1").Value = "imageUrl1"
$ws.Range("G1").Value = "imageUrl2"
$ws.Range("H1").Value = "imageUrl3"

# Data rows - product/brand columns
$ws.Range("B2").Value = "redmi"
$ws.Range("A2").Value = "head phone"
$ws.Range("A3").Value = "earpod"
$ws.Range("A4").Value = "bluetooth"
$ws.Range("B3").Value = "samsunf"
$ws.Range("B4").Value = "zebronics"

# Prices
$ws.Range("C2").Value = 200
$ws.Range("D2").Value = 240
$ws.Range("C3").Value = 300
$ws.Range("D3").Value = 340
$ws.Range("C4").Value = 100
$ws.Range("D4").Value = 430

# Image URLs (columns E-H for all data rows)
$ws.Range("E2").Value = $imgUrl
$ws.Range("F2").Value = $imgUrl
$ws.Range("G2").Value = $imgUrl
$ws.Range("H2").Value = $imgUrl
$ws.Range("E3").Value = $imgUrl
$ws.Range("F3").Value = $imgUrl
$ws.Range("G3").Value = $imgUrl
$ws.Range("H3").Value = $imgUrl
$ws.Range("E4").Value = $imgUrl
$ws.Range("F4").Value = $imgUrl
$ws.Range("G4").Value = $imgUrl
$ws.Range("H4").Value = $imgUrl

# originalUrl header added last (new column E)
$ws.Range("E1").Value = "originalUrl"

# Column widths (values chosen so the saved OOXML "width" attribute matches
# the target 15.5546875 / 19.6640625 / 12.88671875 / 14.21875 / 77 / 75.6640625 / 77 / 76.5546875)
$ws.Columns.Item(1).ColumnWidth = 14.666666666666666
$ws.Columns.Item(2).ColumnWidth = 18.833333333333332
$ws.Columns.Item(3).ColumnWidth = 12.0
$ws.Columns.Item(4).ColumnWidth = 13.333333333333334
$ws.Columns.Item(5).ColumnWidth = 76.16666666666667
$ws.Columns.Item(6).ColumnWidth = 74.83333333333333
$ws.Columns.Item(7).ColumnWidth = 76.16666666666667
$ws.Columns.Item(8).ColumnWidth = 75.66666666666667

# Selection
$ws.Range("B8").Select()
